$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "H 72" record (row 2) entirely, shifting all following rows up by one.
$ws.Rows.Item(2).Delete()
